$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New configuration label for the SMOTE oversampling experiment
$newConfig = "CV + tfidf + ngram(3) + SMOTE"

# Model names, in the same order used for every other configuration block
$models = @("Logistic Regression", "Multinomial Naive Bayes", "Support Vector Machines", "Decision Tree", "Random Forest")

# F1 / Accuracy / Precision / Recall for the 5 new rows (27-31)
$data = @(
    @(88.47, 82.66, 83.97, 93.6),
    @(88.35, 82.66, 84.42, 92.77),
    @(87.98, 80.58, 78.89, 99.79),
    @(80.95, 73.21, 82.51, 79.57),
    @(86.95, 80.41, 83.02, 91.48)
)

$startRow = 27
for ($i = 0; $i -lt $models.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $models[$i]
    $ws.Cells.Item($row, 2).Value = $newConfig
    $ws.Cells.Item($row, 3).Value = $data[$i][0]
    $ws.Cells.Item($row, 4).Value = $data[$i][1]
    $ws.Cells.Item($row, 5).Value = $data[$i][2]
    $ws.Cells.Item($row, 6).Value = $data[$i][3]
}

# The new, longer configuration label widens columns A and B
$ws.Columns.Item(1).ColumnWidth = 21.33
$ws.Columns.Item(2).ColumnWidth = 28

# Scroll the view back to the top and select the last entered cell
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A1").Select() | Out-Null
$ws.Range("F31").Select() | Out-Null
